# Converts an "RRGGBB" hex string into the BGR-packed long that the
# PowerPoint COM object model expects for a ColorFormat.RGB value.
function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-point the three tables (slides 14-16) at the other built-in
#    table style that ships in this deck's table-style list.
# ---------------------------------------------------------------------
$newTableStyleId = "{EA92F32C-B44D-42C3-8269-D94E6719E107}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the presentation's theme away from "Integral" (Red Violet)
#    back to the stock "Office Theme" colour scheme.
# ---------------------------------------------------------------------
$officeColors = @{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$firstSlide = $p.Slides.Item(1)
$colorScheme = $firstSlide.ThemeColorScheme
foreach ($idx in $officeColors.Keys) {
    $colorScheme.Item($idx).RGB = HexToComRgb $officeColors[$idx]
}
